$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRINCIPAL")

$row = 8

$ws.Cells.Item($row, 1).Value = "DF"
$ws.Cells.Item($row, 2).Value = "TESTE00"
$ws.Cells.Item($row, 3).Value = ""
$ws.Cells.Item($row, 4).Value = ""
$ws.Cells.Item($row, 5).Value = ""
$ws.Cells.Item($row, 6).Value = "TESTE"
$ws.Cells.Item($row, 7).Value = "T"
$ws.Cells.Item($row, 8).Value = "T - (T 01/11/25_12H) - DF"
$ws.Cells.Item($row, 9).Value = "01/11/25"
$ws.Cells.Item($row, 10).Value = "12H"
$ws.Cells.Item($row, 11).Value = "14/11/25"
$ws.Cells.Item($row, 12).Value = "DENTRO"
$ws.Cells.Item($row, 13).Value = ""
